$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34
$ws.Range("D34").Value = 44518
$ws.Range("K34").Value = 13000
$ws.Range("L34").Value = 15000
$ws.Range("M34").Value = 13850
$ws.Range("P34").Value = 1065

# Row 35
$ws.Range("D35").Value = 44245
$ws.Range("J35").Value = 400
$ws.Range("K35").Value = 38000
$ws.Range("L35").Value = 40000
$ws.Range("M35").Value = 38850
$ws.Range("P35").Value = 2988

# Row 36
$ws.Range("D36").Value = 44417
$ws.Range("J36").Value = 230
$ws.Range("K36").Value = 13000
$ws.Range("L36").Value = 14000
$ws.Range("M36").Value = 13565
$ws.Range("P36").Value = 1043

# Row 37
$ws.Range("D37").Value = 44445
$ws.Range("J37").Value = 220
$ws.Range("K37").Value = 12000
$ws.Range("L37").Value = 13000
$ws.Range("M37").Value = 12455
$ws.Range("P37").Value = 958

# Row 38
$ws.Range("D38").Value = 44249
$ws.Range("J38").Value = 350
$ws.Range("K38").Value = 38000
$ws.Range("L38").Value = 40000
$ws.Range("M38").Value = 39314
$ws.Range("P38").Value = 3024

# Row 39
$ws.Range("D39").Value = 44270
$ws.Range("J39").Value = 220
$ws.Range("K39").Value = 33000
$ws.Range("L39").Value = 35000
$ws.Range("M39").Value = 34091
$ws.Range("P39").Value = 2622

# Row 40
$ws.Range("D40").Value = 44376
$ws.Range("J40").Value = 580
$ws.Range("K40").Value = 12000
$ws.Range("L40").Value = 14000
$ws.Range("M40").Value = 13103
$ws.Range("P40").Value = 1008

# Row 41
$ws.Range("D41").Value = 44172
$ws.Range("J41").Value = 250
$ws.Range("K41").Value = 27000
$ws.Range("L41").Value = 30000
$ws.Range("M41").Value = 28800
$ws.Range("P41").Value = 2215

# Row 42
$ws.Range("D42").Value = 44242
$ws.Range("J42").Value = 200
$ws.Range("K42").Value = 40000
$ws.Range("L42").Value = 42000
$ws.Range("M42").Value = 41200
$ws.Range("P42").Value = 3169

# Row 43
$ws.Range("D43").Value = 44431
$ws.Range("J43").Value = 260
$ws.Range("K43").Value = 12000
$ws.Range("L43").Value = 13000
$ws.Range("M43").Value = 12462
$ws.Range("P43").Value = 959

# Row 44
$ws.Range("D44").Value = 44301
$ws.Range("J44").Value = 200
$ws.Range("K44").Value = 23000
$ws.Range("M44").Value = 24200
$ws.Range("P44").Value = 1862

# Row 45
$ws.Range("D45").Value = 44284
$ws.Range("J45").Value = 400
$ws.Range("K45").Value = 24000
$ws.Range("L45").Value = 25000
$ws.Range("M45").Value = 24575
$ws.Range("P45").Value = 1890

# Row 46
$ws.Range("D46").Value = 44504
$ws.Range("J46").Value = 230
$ws.Range("K46").Value = 13000
$ws.Range("L46").Value = 15000
$ws.Range("M46").Value = 13870
$ws.Range("P46").Value = 1067

# Row 47
$ws.Range("D47").Value = 44350
$ws.Range("J47").Value = 400
$ws.Range("K47").Value = 23000
$ws.Range("M47").Value = 24150
$ws.Range("P47").Value = 1858

# Row 48
$ws.Range("D48").Value = 44312
$ws.Range("J48").Value = 190
$ws.Range("K48").Value = 24000
$ws.Range("L48").Value = 25000
$ws.Range("M48").Value = 24632
$ws.Range("P48").Value = 1895

# Row 49
$ws.Range("D49").Value = 44382
$ws.Range("J49").Value = 120
$ws.Range("K49").Value = 13000
$ws.Range("L49").Value = 14000
$ws.Range("M49").Value = 13417
$ws.Range("P49").Value = 1032

# Row 50
$ws.Range("D50").Value = 44315
$ws.Range("K50").Value = 25000
$ws.Range("L50").Value = 26000
$ws.Range("M50").Value = 25425
$ws.Range("P50").Value = 1956

# Row 51
$ws.Range("D51").Value = 44510
$ws.Range("J51").Value = 400
$ws.Range("K51").Value = 13000
$ws.Range("L51").Value = 15000
$ws.Range("M51").Value = 13850
$ws.Range("P51").Value = 1065

# Row 52
$ws.Range("D52").Value = 44161
$ws.Range("J52").Value = 330
$ws.Range("K52").Value = 28000
$ws.Range("L52").Value = 30000
$ws.Range("M52").Value = 29394
$ws.Range("P52").Value = 2261

